# Adjust beach trend data to show annual % loss of beaches instead of
# the long-term erosional rate. This inserts two new columns ("years" and
# "trend in % beach lost") after the existing "beach_loss_percent" column,
# pushing the remaining columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new blank columns at I:J - this shifts the former I:N columns
# to K:P, preserving their data/formulas (relative refs get adjusted).
$ws.Columns("I:J").Insert()

# New column headers
$ws.Range("I1").Value = "years"
$ws.Range("J1").Value = "trend in % beach lost"

# "years" column: number of years represented by the long-term erosion data
$ws.Range("I2:I5").Value = 80

# "trend in % beach lost": annual % of beach lost, derived from
# beach_loss_percent (column H) spread out over the number of years (column I)
$ws.Range("J2").Formula = "=-(H2/100)/I2"
$ws.Range("J3:J5").Formula = "=-(H3/100)/I3"

# Re-assert the formulas that got shifted right by the column insert so they
# regroup back into shared formulas spanning rows 3:5, matching rows 2's
# pattern of being the independent "parent" row.
$ws.Range("O3:O5").Formula = "=N3/1000"
$ws.Range("P3:P5").Formula = "=M3*0.001"

# Match the width of the new "years"/"trend" columns to the neighboring
# beach_loss_percent column.
$ws.Columns("I:J").ColumnWidth = $ws.Columns("H").ColumnWidth

# Reflect the author's final cursor position/selection in the sheet view.
$ws.Activate() | Out-Null
$ws.Range("J7").Select() | Out-Null
